$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original text formatting so
# values like "1.003" or "1.000" are not auto-converted to numbers by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.098.45'
$ws.Range('E2').Value = '  -2.83%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.640.77'
$ws.Range('E3').Value = '  -2.57%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.77'
$ws.Range('E5').Value = '  -2.00%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  -0.18%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3898'
$ws.Range('E7').Value = '  -0.93%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3853'
$ws.Range('E8').Value = '  -2.94%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.004'
$ws.Range('E9').Value = '  -0.12%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '49.44'
$ws.Range('E10').Value = '  -4.63%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.351'
$ws.Range('E11').Value = '  -5.12%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08592'
$ws.Range('E12').Value = '  -0.98%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.075'
$ws.Range('E13').Value = '  -3.40%  '

$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.34'
$ws.Range('E14').Value = '  -7.34%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001286'
$ws.Range('E15').Value = '  -2.69%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.450'
$ws.Range('E16').Value = '  -4.44%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.648.89'
$ws.Range('E17').Value = '  +16.54%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '94.89'
$ws.Range('E18').Value = '  +0.70%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06906'
$ws.Range('E19').Value = '  -2.75%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.30'
$ws.Range('E20').Value = '  +0.70%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.879'
$ws.Range('E21').Value = '  -3.68%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  -0.20%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.52'
$ws.Range('E23').Value = '  -4.22%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '24.118.58'
$ws.Range('E24').Value = '  -2.71%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.422'
$ws.Range('E25').Value = '  +2.25%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.806'
$ws.Range('E26').Value = '  +1.09%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.33'
$ws.Range('E27').Value = '  -6.25%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '157.54'
$ws.Range('E28').Value = '  -3.03%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.447'
$ws.Range('E29').Value = '  +7.64%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.365'
$ws.Range('E30').Value = '  -6.50%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '140.11'
$ws.Range('E31').Value = '  -7.00%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.403'
$ws.Range('E32').Value = '  -7.44%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.824.60'
$ws.Range('E33').Value = '  -2.16%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.959'
$ws.Range('E34').Value = '  +0.23%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08080'
$ws.Range('E35').Value = '  -4.56%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02889'
$ws.Range('E36').Value = '  -6.13%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2674'
$ws.Range('E37').Value = '  -4.80%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9456'
$ws.Range('E38').Value = '  -6.90%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.09189'
$ws.Range('E39').Value = '  -4.06%  '

$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.04'
$ws.Range('E40').Value = '  -4.42%  '

$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.455'
$ws.Range('E41').Value = '  -1.22%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7505'
$ws.Range('E42').Value = '  -5.71%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.94'
$ws.Range('E43').Value = '  -5.60%  '

$ws.Range('E44').Value = '  -4.06%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6878'
$ws.Range('E45').Value = '  -4.03%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.447'
$ws.Range('E46').Value = '  -5.33%  '

$ws.Range('E47').Value = '  -2.54%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.000'
$ws.Range('E48').Value = '  -0.25%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08348'
$ws.Range('E49').Value = '  -4.65%  '

$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.89'
$ws.Range('E50').Value = '  -4.16%  '

$ws.Range('B51').Value = 'Flow'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.252'
$ws.Range('E51').Value = '  -6.87%  '
